$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap data rows 55 and 56 (keep column A index in place)
$ws.Range("B55").Value = 6830656
$ws.Range("B56").Value = 6830657
$ws.Range("F55").Value = "Mazatlan FC Women"
$ws.Range("F56").Value = "Chivas Guadalajara Women"
$ws.Range("G55").Value = "Club Necaxa Women"
$ws.Range("G56").Value = "Monterrey Women"
$ws.Range("H55").Value = 0
$ws.Range("H56").Value = 1
$ws.Range("I55").Value = 2
$ws.Range("I56").Value = 0
$ws.Range("J55").Value = "A"
$ws.Range("J56").Value = "H"
$ws.Range("K55").Value = 3.5
$ws.Range("K56").Value = 2.375
$ws.Range("L55").Value = 3.5
$ws.Range("L56").Value = 3.4
$ws.Range("M55").Value = 1.833
$ws.Range("M56").Value = 2.5
$ws.Range("N55").Value = 4.333
$ws.Range("N56").Value = 2.15
$ws.Range("O55").Value = 3.75
$ws.Range("O56").Value = 3.5
$ws.Range("P55").Value = 1.727
$ws.Range("P56").Value = 3
$ws.Range("Q55").Value = 0.75
$ws.Range("Q56").Value = -0.25
$ws.Range("R55").Value = 1.85
$ws.Range("R56").Value = 1.875
$ws.Range("S55").Value = 1.95
$ws.Range("S56").Value = 1.925
$ws.Range("U55").Value = 1.975
$ws.Range("U56").Value = 1.75
$ws.Range("V55").Value = 1.825
$ws.Range("V56").Value = 1.95
$ws.Range("W55").Value = -1
$ws.Range("W56").Value = 1.15
$ws.Range("Y55").Value = 0.7270000000000001
$ws.Range("Y56").Value = -1
$ws.Range("Z55").Value = -1
$ws.Range("Z56").Value = 0.875
$ws.Range("AA55").Value = 0.95
$ws.Range("AA56").Value = -1
$ws.Range("AC55").Value = 0.825
$ws.Range("AC56").Value = 0.95

# Swap data rows 109 and 110 (keep column A index in place)
$ws.Range("B109").Value = 6830712
$ws.Range("B110").Value = 6830711
$ws.Range("F109").Value = "Tigres UANL Women"
$ws.Range("F110").Value = "Atletico San Luis Women"
$ws.Range("G109").Value = "Unam Pumas Women"
$ws.Range("G110").Value = "Tijuana Women"
$ws.Range("H109").Value = 3
$ws.Range("H110").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("I110").Value = 2
$ws.Range("J109").Value = "H"
$ws.Range("J110").Value = "A"
$ws.Range("K109").Value = 1.181
$ws.Range("K110").Value = 3
$ws.Range("L109").Value = 6
$ws.Range("L110").Value = 3.6
$ws.Range("M109").Value = 10
$ws.Range("M110").Value = 2
$ws.Range("N109").Value = 1.1
$ws.Range("N110").Value = 4
$ws.Range("O109").Value = 9
$ws.Range("O110").Value = 3.8
$ws.Range("P109").Value = 19
$ws.Range("P110").Value = 1.666
$ws.Range("Q109").Value = -2.5
$ws.Range("Q110").Value = 0.75
$ws.Range("R109").Value = 1.8
$ws.Range("R110").Value = 1.925
$ws.Range("S109").Value = 2
$ws.Range("S110").Value = 1.875
$ws.Range("T109").Value = 3.75
$ws.Range("T110").Value = 3
$ws.Range("W109").Value = 0.1000000000000001
$ws.Range("W110").Value = -1
$ws.Range("Y109").Value = -1
$ws.Range("Y110").Value = 0.6659999999999999
$ws.Range("Z109").Value = 0.8
$ws.Range("Z110").Value = -1
$ws.Range("AA109").Value = -1
$ws.Range("AA110").Value = 0.875

# Swap data rows 213 and 214 (keep column A index in place)
$ws.Range("B213").Value = 7645771
$ws.Range("B214").Value = 7645770
$ws.Range("F213").Value = "Toluca Women"
$ws.Range("F214").Value = "Atletico San Luis Women"
$ws.Range("G213").Value = "Mazatlan FC Women"
$ws.Range("G214").Value = "Unam Pumas Women"
$ws.Range("H213").Value = 6
$ws.Range("H214").Value = 0
$ws.Range("I213").Value = 0
$ws.Range("I214").Value = 3
$ws.Range("J213").Value = "H"
$ws.Range("J214").Value = "A"
$ws.Range("K213").Value = 1.142
$ws.Range("K214").Value = 4.333
$ws.Range("L213").Value = 7
$ws.Range("L214").Value = 4
$ws.Range("M213").Value = 13
$ws.Range("M214").Value = 1.666
$ws.Range("N213").Value = 1.071
$ws.Range("N214").Value = 3.5
$ws.Range("O213").Value = 10
$ws.Range("O214").Value = 4
$ws.Range("P213").Value = 23
$ws.Range("P214").Value = 1.833
$ws.Range("Q213").Value = -2.5
$ws.Range("Q214").Value = 0.5
$ws.Range("R213").Value = 1.75
$ws.Range("R214").Value = 1.95
$ws.Range("S213").Value = 2.05
$ws.Range("S214").Value = 1.85
$ws.Range("T213").Value = 3.75
$ws.Range("T214").Value = 3.25
$ws.Range("U213").Value = 1.9
$ws.Range("U214").Value = 1.95
$ws.Range("V213").Value = 1.9
$ws.Range("V214").Value = 1.75
$ws.Range("W213").Value = 0.07099999999999995
$ws.Range("W214").Value = -1
$ws.Range("Y213").Value = -1
$ws.Range("Y214").Value = 0.833
$ws.Range("Z213").Value = 0.75
$ws.Range("Z214").Value = -1
$ws.Range("AA213").Value = -1
$ws.Range("AA214").Value = 0.8500000000000001
$ws.Range("AB213").Value = 0.8999999999999999
$ws.Range("AB214").Value = -0.5
$ws.Range("AC213").Value = -1
$ws.Range("AC214").Value = 0.375

# Swap data rows 215 and 216 (keep column A index in place)
$ws.Range("B215").Value = 7645772
$ws.Range("B216").Value = 7645707
$ws.Range("F215").Value = "Pachuca Women"
$ws.Range("F216").Value = "Tigres UANL Women"
$ws.Range("G215").Value = "Queretaro Women"
$ws.Range("G216").Value = "Tijuana Women"
$ws.Range("H215").Value = 4
$ws.Range("H216").Value = 2
$ws.Range("I215").Value = 1
$ws.Range("I216").Value = 0
$ws.Range("K215").Value = 1.25
$ws.Range("K216").Value = 1.125
$ws.Range("L215").Value = 5.5
$ws.Range("L216").Value = 7.5
$ws.Range("M215").Value = 7.5
$ws.Range("M216").Value = 15
$ws.Range("N215").Value = 1.285
$ws.Range("N216").Value = 1.166
$ws.Range("O215").Value = 5.5
$ws.Range("O216").Value = 7
$ws.Range("P215").Value = 6.5
$ws.Range("P216").Value = 15
$ws.Range("Q215").Value = -1.75
$ws.Range("Q216").Value = -2.25
$ws.Range("R215").Value = 1.975
$ws.Range("R216").Value = 1.825
$ws.Range("S215").Value = 1.825
$ws.Range("S216").Value = 1.975
$ws.Range("T215").Value = 3.25
$ws.Range("T216").Value = 3.75
$ws.Range("U215").Value = 1.8
$ws.Range("U216").Value = 1.85
$ws.Range("V215").Value = 2
$ws.Range("V216").Value = 1.95
$ws.Range("W215").Value = 0.2849999999999999
$ws.Range("W216").Value = 0.1659999999999999
$ws.Range("Z215").Value = 0.9750000000000001
$ws.Range("Z216").Value = -0.5
$ws.Range("AA215").Value = -1
$ws.Range("AA216").Value = 0.4875
$ws.Range("AB215").Value = 0.8
$ws.Range("AB216").Value = -1
$ws.Range("AC215").Value = -1
$ws.Range("AC216").Value = 0.95

# Swap data rows 248 and 249 (keep column A index in place)
$ws.Range("B248").Value = 7645796
$ws.Range("B249").Value = 7645719
$ws.Range("F248").Value = "Chivas Guadalajara Women"
$ws.Range("F249").Value = "Monterrey Women"
$ws.Range("G248").Value = "Club Necaxa Women"
$ws.Range("G249").Value = "Santos Laguna Women"
$ws.Range("H248").Value = 4
$ws.Range("H249").Value = 6
$ws.Range("K248").Value = 1.055
$ws.Range("K249").Value = 1.025
$ws.Range("L248").Value = 10
$ws.Range("L249").Value = 15
$ws.Range("M248").Value = 21
$ws.Range("M249").Value = 34
$ws.Range("N248").Value = 1.062
$ws.Range("N249").Value = 1.01
$ws.Range("O248").Value = 11
$ws.Range("O249").Value = 34
$ws.Range("P248").Value = 29
$ws.Range("P249").Value = 67
$ws.Range("Q248").Value = -3
$ws.Range("Q249").Value = -4.75
$ws.Range("R248").Value = 1.85
$ws.Range("R249").Value = 1.775
$ws.Range("S248").Value = 1.95
$ws.Range("S249").Value = 1.925
$ws.Range("T248").Value = 4
$ws.Range("T249").Value = 5.75
$ws.Range("U248").Value = 1.8
$ws.Range("U249").Value = 1.85
$ws.Range("V248").Value = 2
$ws.Range("V249").Value = 1.95
$ws.Range("W248").Value = 0.06200000000000006
$ws.Range("W249").Value = 0.01000000000000001
$ws.Range("Z248").Value = 0.8500000000000001
$ws.Range("Z249").Value = 0.7749999999999999
$ws.Range("AB248").Value = 0
$ws.Range("AB249").Value = 0.425
$ws.Range("AC248").Value = -0
$ws.Range("AC249").Value = -0.5

# Swap data rows 263 and 265 (keep column A index in place)
$ws.Range("B263").Value = 7645804
$ws.Range("B265").Value = 7645807
$ws.Range("F263").Value = "Puebla Women"
$ws.Range("F265").Value = "Club Necaxa Women"
$ws.Range("G263").Value = "Mazatlan FC Women"
$ws.Range("G265").Value = "Leon Women"
$ws.Range("H263").Value = 1
$ws.Range("H265").Value = 2
$ws.Range("I263").Value = 2
$ws.Range("I265").Value = 1
$ws.Range("J263").Value = "A"
$ws.Range("J265").Value = "H"
$ws.Range("K263").Value = 1.666
$ws.Range("K265").Value = 4.333
$ws.Range("M263").Value = 3.8
$ws.Range("M265").Value = 1.571
$ws.Range("N263").Value = 1.333
$ws.Range("N265").Value = 7
$ws.Range("O263").Value = 4.5
$ws.Range("O265").Value = 4.2
$ws.Range("P263").Value = 7.5
$ws.Range("P265").Value = 1.363
$ws.Range("Q263").Value = -1.5
$ws.Range("Q265").Value = 1.5
$ws.Range("R263").Value = 1.925
$ws.Range("R265").Value = 1.75
$ws.Range("S263").Value = 1.875
$ws.Range("S265").Value = 1.95
$ws.Range("T263").Value = 3
$ws.Range("T265").Value = 2.75
$ws.Range("U263").Value = 1.75
$ws.Range("U265").Value = 1.8
$ws.Range("V263").Value = 2.05
$ws.Range("V265").Value = 2
$ws.Range("W263").Value = -1
$ws.Range("W265").Value = 6
$ws.Range("Y263").Value = 6.5
$ws.Range("Y265").Value = -1
$ws.Range("Z263").Value = -1
$ws.Range("Z265").Value = 0.75
$ws.Range("AA263").Value = 0.875
$ws.Range("AA265").Value = -1
$ws.Range("AB263").Value = 0
$ws.Range("AB265").Value = 0.4
$ws.Range("AC263").Value = -0
$ws.Range("AC265").Value = -0.5

